# Update OPIS prompt and mappings.xlsx for improved terminal extraction
#
# This script updates the "OPIS_Terminal_Mappings" worksheet:
#   - Removes the old "ST. LOUIS, MO" -> "St. Louis, MO" row (old row 19),
#     which shifts all the following rows up by one.
#   - Re-applies the compact (non-wrapped, 12.8pt) formatting to the
#     "East St. Louis, IL" -> "IL East St. Louis" row, which is now row 32.
#   - Appends two new mapping rows at the bottom of the table:
#       33: "ST. LOUIS TERMINAL" -> "St. Louis, MO"
#       34: "FOB St. Louis"      -> "St. Louis, MO"
#   - Updates the sheet view so the window is scrolled down a bit further
#     and the new blank row below the table is selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OPIS_Terminal_Mappings")
$ws.Activate()

# Remove the obsolete "ST. LOUIS, MO" row; everything below shifts up by 1.
$ws.Rows.Item(19).Delete()

# The "East St. Louis, IL" / "IL East St. Louis" row is now row 32; give it
# the compact, non-wrapped formatting used by the new trailer rows.
$ws.Rows.Item(32).RowHeight = 12.8
$ws.Range("A32:B32").WrapText = $False

# New row 33: "ST. LOUIS TERMINAL" -> "St. Louis, MO"
$ws.Range("A33").Value = "ST. LOUIS TERMINAL"
$ws.Range("B33").Value = "St. Louis, MO"
$ws.Rows.Item(33).RowHeight = 12.8
$ws.Range("A33:B33").WrapText = $False

# New row 34: "FOB St. Louis" -> "St. Louis, MO"
$ws.Range("A34").Value = "FOB St. Louis"
$ws.Range("B34").Value = "St. Louis, MO"
$ws.Rows.Item(34).RowHeight = 12.8
$ws.Range("A34:B34").WrapText = $False

# Update the view: scroll so row 19 is at the top and select the first
# empty row beneath the table (A35), matching the new window state.
$ws.Range("A19").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A35").Select()
